$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.675.94'
$ws.Range('E2').Value = '  +0.45%  '
$ws.Range('D3').Value = '1.849.70'
$ws.Range('E3').Value = '  +0.49%  '
$ws.Range('D4').Value = '1.002'
$ws.Range('E4').Value = '  +0.25%  '
$ws.Range('D5').Value = '312.28'
$ws.Range('E5').Value = '  -0.19%  '
$ws.Range('E6').Value = '  +0.14%  '
$ws.Range('D7').Value = '0.4255'
$ws.Range('E7').Value = '  +0.55%  '
$ws.Range('D8').Value = '0.3634'
$ws.Range('E8').Value = '  +0.11%  '
$ws.Range('D9').Value = '44.68'
$ws.Range('E9').Value = '  +2.40%  '
$ws.Range('D10').Value = '0.07279'
$ws.Range('E10').Value = '  +1.11%  '
$ws.Range('D11').Value = '0.8723'
$ws.Range('E11').Value = '  -2.91%  '
$ws.Range('D12').Value = '20.55'
$ws.Range('E12').Value = '  -0.25%  '
$ws.Range('D13').Value = '1.904.82'
$ws.Range('E13').Value = '  +3.71%  '
$ws.Range('D14').Value = '5.307'
$ws.Range('E14').Value = '  +0.00%  '
$ws.Range('D15').Value = '6.500'
$ws.Range('E15').Value = '  -1.11%  '
$ws.Range('D16').Value = '0.06894'
$ws.Range('E16').Value = '  +1.40%  '
$ws.Range('D17').Value = '1.003'
$ws.Range('E17').Value = '  +0.19%  '
$ws.Range('D18').Value = '79.70'
$ws.Range('E18').Value = '  +3.49%  '
$ws.Range('D19').Value = '0.000009028'
$ws.Range('E19').Value = '  +1.51%  '
$ws.Range('D20').Value = '1.001'
$ws.Range('E20').Value = '  +0.14%  '
$ws.Range('E21').Value = '  +0.24%  '
$ws.Range('D22').Value = '27.696.26'
$ws.Range('E22').Value = '  +0.62%  '
$ws.Range('D23').Value = '4.964'
$ws.Range('D24').Value = '10.37'
$ws.Range('E24').Value = '  -3.37%  '
$ws.Range('D25').Value = '2.115.97'
$ws.Range('E25').Value = '  +3.11%  '
$ws.Range('D26').Value = '1.961'
$ws.Range('E26').Value = '  -4.07%  '
$ws.Range('D27').Value = '153.52'
$ws.Range('E27').Value = '  +1.58%  '
$ws.Range('E28').Value = '  +3.30%  '
$ws.Range('D29').Value = '121.93'
$ws.Range('E29').Value = '  +10.28%  '
$ws.Range('D30').Value = '5.256'
$ws.Range('E30').Value = '  -0.17%  '
$ws.Range('D31').Value = '1.864'
$ws.Range('E31').Value = '  +10.87%  '
$ws.Range('D32').Value = '0.08890'
$ws.Range('E32').Value = '  +0.42%  '
$ws.Range('D33').Value = '0.7599'
$ws.Range('E33').Value = '  -1.62%  '
$ws.Range('E34').Value = '  +4.29%  '
$ws.Range('D35').Value = '4.511'
$ws.Range('D36').Value = '1.095'
$ws.Range('E36').Value = '  +1.27%  '
$ws.Range('D37').Value = '0.05370'
$ws.Range('E37').Value = '  -0.20%  '
$ws.Range('D38').Value = '1.090'
$ws.Range('E38').Value = '  -0.41%  '
$ws.Range('D39').Value = '0.01927'
$ws.Range('E39').Value = '  +0.42%  '
$ws.Range('D40').Value = '2.819'
$ws.Range('E40').Value = '  -4.16%  '
$ws.Range('D41').Value = '0.5044'
$ws.Range('E41').Value = '  +0.22%  '
$ws.Range('D42').Value = '0.1643'
$ws.Range('E42').Value = '  +0.61%  '
$ws.Range('D43').Value = '6.746'
$ws.Range('E43').Value = '  -0.32%  '
$ws.Range('D44').Value = '8.328'
$ws.Range('E44').Value = '  +1.48%  '
$ws.Range('D45').Value = '0.06537'
$ws.Range('E45').Value = '  -1.10%  '
$ws.Range('D46').Value = '10.34'
$ws.Range('E46').Value = '  +1.73%  '
$ws.Range('D47').Value = '104.80'
$ws.Range('E47').Value = '  -0.62%  '
$ws.Range('D48').Value = '0.4640'
$ws.Range('E48').Value = '  -1.60%  '
$ws.Range('D49').Value = '0.9993'
$ws.Range('E49').Value = '  +0.03%  '
$ws.Range('D50').Value = '1.617'
$ws.Range('E50').Value = '  -1.38%  '
$ws.Range('D51').Value = '64.37'
$ws.Range('E51').Value = '  +0.01%  '
